$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pass 1: ids (column C) - proc_1..proc_10 (proc_1..proc_9 already exist, proc_10 is new)
for ($i = 1; $i -le 10; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = "proc_$i"
}

# Pass 2: names (column D) - process_1..process_10 (all new)
for ($i = 1; $i -le 10; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 4).Value = "process_$i"
}

# Pass 3: remaining columns, row by row, in column order
for ($i = 1; $i -le 10; $i++) {
    $row = $i + 1

    $ws.Cells.Item($row, 2).Value = "env"                 # B: env
    $ws.Cells.Item($row, 5).Value = 14 + $i               # E: proc_time
    $ws.Cells.Item($row, 7).Value = 2 * $i - 1            # G: operators
    $ws.Cells.Item($row, 8).Value = "available"          # H: operating_status
    $ws.Cells.Item($row, 9).Value = "[]"                 # I: upstream_processes
    $ws.Cells.Item($row, 10).Value = "[]"                # J: downstream_processes
    $ws.Cells.Item($row, 11).Value = "[]"                # K: sub_processes
    $ws.Cells.Item($row, 12).Value = "[]"                # L: skills
    $ws.Cells.Item($row, 13).Value = "{}"                # M: input_products
    $ws.Cells.Item($row, 14).Value = "{}"                # N: output_products
    $ws.Cells.Item($row, 15).Value = "{}"                # O: resources
}

# Columns F (operating_cost header), G (operators header) and H (operating_status
# header / "available" values) grew to fit their header captions once the data
# below was populated, exactly like Excel's "best fit" column sizing does on import.
$ws.Columns.Item(6).ColumnWidth = 13.45
$ws.Columns.Item(7).ColumnWidth = 8.65
$ws.Columns.Item(8).ColumnWidth = 15.16

$ws.Range("P2").Select()
